$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mobility")

# Row 24 holds the main_min_density parameter.
$ws.Range("B24").Value = 800000
$ws.Range("C24").Value = "Minimum net density to consider a link as being a main track (ton-km/ton = ton)."
